$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing headers (e.g. H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 and IF columns with their data values
$data = @(
    @(5, 6),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(5, 6),
    @(6, 6),
    @(5, 5),
    @(5, 6),
    @(6, 7),
    @(5, 5),
    @(9, 9),
    @(5, 6),
    @(6, 7),
    @(5, 5),
    @(9, 9),
    @(5, 5),
    @(5, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
